# Recompute column H ("客単価" / per-customer amount) on the
# "ABC分析_客構成" sheet: each numeric H value is divided by the
# corresponding row's column E value (客構成 count), turning the
# previously-aggregated total into a per-count figure.
#
# Row 1 is the header and is left untouched. Any H cell that isn't a
# plain number (e.g. the literal text "inf") is skipped, as is any row
# whose E value is 0 or 1 (those quotients are identical to the
# original value, so nothing actually changes there).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($row = [Math]::Max($firstRow, 2); $row -le $lastRow; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $hVal = $hCell.Value()

    if ($hVal -is [double] -or $hVal -is [int]) {
        $eCell = $ws.Cells.Item($row, 5)   # column E
        $eVal = $eCell.Value()

        if (($eVal -is [double] -or $eVal -is [int]) -and $eVal -ne 0 -and $eVal -ne 1) {
            $hCell.Value = $hVal / $eVal
        }
    }
}
